# Auto-generated Excel COM-interop script
# Applies the scraped-data refresh described in the commit:
# "Horarios actualizados Linea 141 - 439"
# It updates the 'last updated' / 'total rows' headers and the
# Hora_Scrap/Hora_Llegada/Linea/Minutos/Parada data table on each
# of the 3 worksheets (LP1912, LP1912-215, 6203-6173).

$wb = $excel.ActiveWorkbook

# ----- Sheet: LP1912 -----
$ws = $wb.Worksheets.Item('LP1912')

$ws.Range('A2').Value = 'Última actualización: 12:24:14'
$ws.Range('A3').Value = 'Total filas: 245'

$rows = @(
    ,@(65, '07:14:27', '07:37', '23_HERNANDEZ', 23, 'LP1912')
    ,@(66, '06:02:16', '07:37', '27_EL RETIRO', 95, 'LP1912')
    ,@(111, '08:47:19', '09:23', '16_SANTA ANA', 36, 'LP1912')
    ,@(112, '07:44:08', '09:23', '17_ROMERO', 99, 'LP1912')
    ,@(113, '07:57:27', '09:23', '11_ETCHEVERRY', 86, 'LP1912')
    ,@(188, '11:15:53', '11:58', '225_GOMEZ', 43, 'LP1912')
    ,@(189, '11:58:46', '11:58', '17_ROMERO', 0, 'LP1912')
    ,@(197, '10:11:11', '12:07', '16_P MOR-SANTA ANA', 116, 'LP1912')
    ,@(198, '10:50:37', '12:07', '10_OLMOS', 77, 'LP1912')
    ,@(199, '10:11:11', '12:07', '14_ABASTO', 116, 'LP1912')
    ,@(206, '10:50:37', '12:21', '26_HERNANDEZ', 91, 'LP1912')
    ,@(207, '10:50:37', '12:21', '14_ABASTO', 91, 'LP1912')
    ,@(209, '12:24:14', '12:24', '16_SANTA ANA', 0, 'LP1912')
    ,@(210, '12:24:14', '12:24', '17_ROMERO', 0, 'LP1912')
    ,@(211, '12:24:14', '12:27', '10_OLMOS', 3, 'LP1912')
    ,@(212, '11:15:53', '12:34', '11_ETCHEVERRY', 79, 'LP1912')
    ,@(213, '11:58:46', '12:34', '23_HERNANDEZ', 36, 'LP1912')
    ,@(214, '12:24:14', '12:35', '23_HERNANDEZ', 11, 'LP1912')
    ,@(215, '12:24:14', '12:35', '11_ETCHEVERRY', 11, 'LP1912')
    ,@(216, '10:50:37', '12:36', '27_EL RETIRO', 106, 'LP1912')
    ,@(217, '11:43:19', '12:37', '27_EL RETIRO', 54, 'LP1912')
    ,@(218, '11:43:19', '12:37', '23_HERNANDEZ', 54, 'LP1912')
    ,@(219, '10:50:37', '12:38', '17_179 Y 38', 108, 'LP1912')
    ,@(220, '11:15:53', '12:40', '10_OLMOS', 85, 'LP1912')
    ,@(221, '11:43:19', '12:41', '10_OLMOS', 58, 'LP1912')
    ,@(222, '11:15:53', '12:46', '17_ROMERO', 91, 'LP1912')
    ,@(223, '12:24:14', '12:46', '16_SANTA ANA', 22, 'LP1912')
    ,@(224, '11:15:53', '12:48', '11_ETCHEVERRY', 93, 'LP1912')
    ,@(225, '12:24:14', '12:49', '11_ETCHEVERRY', 25, 'LP1912')
    ,@(226, '11:15:53', '13:02', '15_ABASTO', 107, 'LP1912')
    ,@(227, '12:24:14', '13:03', '14_ABASTO', 39, 'LP1912')
    ,@(228, '11:15:53', '13:06', '16_P MOR-SANTA ANA', 111, 'LP1912')
    ,@(229, '12:24:14', '13:07', '16_P MOR-SANTA ANA', 43, 'LP1912')
    ,@(230, '11:15:53', '13:13', '215D_EL PATO', 118, 'LP1912')
    ,@(231, '11:43:19', '13:14', '215D_EL PATO', 91, 'LP1912')
    ,@(232, '11:43:19', '13:14', '17_ROMERO', 91, 'LP1912')
    ,@(233, '11:43:19', '13:19', '10_OLMOS', 96, 'LP1912')
    ,@(234, '11:43:19', '13:21', '26_HERNANDEZ', 98, 'LP1912')
    ,@(235, '11:43:19', '13:26', '14_ABASTO', 103, 'LP1912')
    ,@(236, '11:43:19', '13:26', '15_ABASTO', 103, 'LP1912')
    ,@(237, '12:24:14', '13:27', '14_ABASTO', 63, 'LP1912')
    ,@(238, '11:58:46', '13:29', '17_ROMERO', 91, 'LP1912')
    ,@(239, '11:43:19', '13:37', '11_ETCHEVERRY', 114, 'LP1912')
    ,@(240, '12:24:14', '13:40', '23_HERNANDEZ', 76, 'LP1912')
    ,@(241, '11:58:46', '13:46', '17_ROMERO', 108, 'LP1912')
    ,@(242, '12:24:14', '13:47', '17_ROMERO', 83, 'LP1912')
    ,@(243, '11:58:46', '13:50', '215A_EL PATO', 112, 'LP1912')
    ,@(244, '12:24:14', '13:51', '215A_EL PATO', 87, 'LP1912')
    ,@(245, '11:58:46', '13:56', '16_P MOR-167 Y 521', 118, 'LP1912')
    ,@(246, '11:58:46', '13:56', '225_GOMEZ', 118, 'LP1912')
    ,@(247, '12:24:14', '13:57', '16_P MOR-167 Y 521', 93, 'LP1912')
    ,@(248, '12:24:14', '14:17', '27_EL RETIRO', 113, 'LP1912')
    ,@(249, '12:24:14', '14:20', '215C_EL PATO', 116, 'LP1912')
    ,@(250, '12:24:14', '14:21', '26_HERNANDEZ', 117, 'LP1912')
)

foreach ($row in $rows) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
    $ws.Cells.Item($r, 5).Value = $row[5]
}

# ----- Sheet: LP1912-215 -----
$ws = $wb.Worksheets.Item('LP1912-215')

$ws.Range('A2').Value = 'Última actualización: 12:24:14'
$ws.Range('A3').Value = 'Total filas: 29'

$rows = @(
    ,@(33, '12:24:14', '13:51', '215A_EL PATO', 87, 'LP1912')
    ,@(34, '12:24:14', '14:20', '215C_EL PATO', 116, 'LP1912')
)

foreach ($row in $rows) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
    $ws.Cells.Item($r, 5).Value = $row[5]
}

# ----- Sheet: 6203-6173 -----
$ws = $wb.Worksheets.Item('6203-6173')

$ws.Range('A2').Value = 'Última actualización: 12:24:14'
$ws.Range('A3').Value = 'Total filas: 36'

$rows = @(
    ,@(41, '12:24:14', '14:09', '215A_LA PLATA', 105, 'L6173')
)

foreach ($row in $rows) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
    $ws.Cells.Item($r, 5).Value = $row[5]
}

